{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,style,text\");\nawait context.sync();\n\n// Remove the errant trailing semicolon from the \"_h2o_keep_element;\"\n// marker text in every \"Node End\" paragraph (export/node.html).\nfor (const p of paragraphs.items) {\n  if (p.style === \"Node End\" && p.text === \"_h2o_keep_element;\") {\n    const results = p.search(\"_h2o_keep_element;\", { matchCase: true });\n    results.load(\"items\");\n    await context.sync();\n    for (const r of results.items) {\n      r.insertText(\"_h2o_keep_element\", Word.InsertLocation.replace);\n    }\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the errant trailing semicolon from the \"_h2o_keep_element;\"\n# marker text in every \"Node End\" paragraph (export/node.html).\nforeach ($p in $d.Paragraphs) {\n    if ($p.Style.NameLocal -eq \"Node End\") {\n        $t = $p.Range.Text.TrimEnd([char]13)\n        if ($t -eq \"_h2o_keep_element;\") {\n            $p.Range.Text = $t.Substring(0, $t.Length - 1)\n        }\n    }\n}\n"}
